$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so values are stored verbatim
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "30.414.58"
$ws.Cells.Item(3, 4).Value = "2.017.35"
$ws.Cells.Item(3, 5).Value = "  +5.91%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "324.89"
$ws.Cells.Item(5, 5).Value = "  +1.73%  "
$ws.Cells.Item(6, 5).Value = "  +0.00%  "
$ws.Cells.Item(7, 4).Value = "0.5136"
$ws.Cells.Item(7, 5).Value = "  +2.17%  "
$ws.Cells.Item(8, 4).Value = "0.4221"
$ws.Cells.Item(8, 5).Value = "  +4.45%  "
$ws.Cells.Item(9, 4).Value = "0.08715"
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(10, 4).Value = "1.136"
$ws.Cells.Item(10, 5).Value = "  +3.84%  "
$ws.Cells.Item(11, 2).Value = "Solana"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(11, 4).Value = "24.90"
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "2.011.46"
$ws.Cells.Item(12, 5).Value = "  +5.66%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "6.608"
$ws.Cells.Item(13, 5).Value = "  +4.03%  "
$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).Value = "7.483"
$ws.Cells.Item(14, 5).Value = "  +4.37%  "
$ws.Cells.Item(15, 2).Value = "BinanceUSD"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(15, 4).Value = "1.003"
$ws.Cells.Item(15, 5).Value = "  +0.09%  "
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).Value = "94.51"
$ws.Cells.Item(16, 5).Value = "  +2.93%  "
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).Value = "0.00001114"
$ws.Cells.Item(17, 5).Value = "  +2.08%  "
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(18, 4).Value = "0.06529"
$ws.Cells.Item(18, 5).Value = "  +0.63%  "
$ws.Cells.Item(19, 2).Value = "Avalanche"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(19, 4).Value = "18.98"
$ws.Cells.Item(19, 5).Value = "  +5.49%  "
$ws.Cells.Item(20, 2).Value = "Dai"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(20, 4).Value = "1.001"
$ws.Cells.Item(20, 5).Value = "  +0.00%  "
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).Value = "6.209"
$ws.Cells.Item(21, 5).Value = "  +4.68%  "
$ws.Cells.Item(22, 2).Value = "WrappedBTC"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(22, 4).Value = "30.481.02"
$ws.Cells.Item(22, 5).Value = "  +1.67%  "
$ws.Cells.Item(23, 2).Value = "Cosmos"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(23, 4).Value = "11.87"
$ws.Cells.Item(23, 5).Value = "  +5.74%  "
$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(24, 4).Value = "2.226"
$ws.Cells.Item(24, 5).Value = "  +1.65%  "
$ws.Cells.Item(25, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(25, 4).Value = "2.254.33"
$ws.Cells.Item(25, 5).Value = "  +6.14%  "
$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).Value = "22.47"
$ws.Cells.Item(26, 5).Value = "  +2.12%  "
$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).Value = "162.72"
$ws.Cells.Item(27, 5).Value = "  +0.97%  "
$ws.Cells.Item(28, 2).Value = "LidoDAOToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(28, 4).Value = "2.427"
$ws.Cells.Item(28, 5).Value = "  +7.25%  "
$ws.Cells.Item(29, 2).Value = "BitcoinCash"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(29, 4).Value = "131.55"
$ws.Cells.Item(29, 5).Value = "  +2.35%  "
$ws.Cells.Item(30, 2).Value = "ImmutableX"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(30, 4).Value = "1.141"
$ws.Cells.Item(30, 5).Value = "  +1.86%  "
$ws.Cells.Item(31, 2).Value = "Stellar"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(31, 4).Value = "0.1053"
$ws.Cells.Item(31, 5).Value = "  +2.01%  "
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).Value = "6.071"
$ws.Cells.Item(32, 5).Value = "  +3.02%  "
$ws.Cells.Item(33, 2).Value = "HuobiToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(33, 4).Value = "3.831"
$ws.Cells.Item(33, 5).Value = "  +0.91%  "
$ws.Cells.Item(34, 2).Value = "ARBITRUM"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(34, 4).Value = "1.371"
$ws.Cells.Item(34, 5).Value = "  +15.39%  "
$ws.Cells.Item(35, 2).Value = "VeChain"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(35, 4).Value = "0.02527"
$ws.Cells.Item(35, 5).Value = "  +4.03%  "
$ws.Cells.Item(36, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(36, 4).Value = "5.486"
$ws.Cells.Item(36, 5).Value = "  +2.35%  "
$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).Value = "0.06666"
$ws.Cells.Item(37, 5).Value = "  +5.38%  "
$ws.Cells.Item(38, 2).Value = "Aptos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(38, 4).Value = "12.31"
$ws.Cells.Item(38, 5).Value = "  +9.02%  "
$ws.Cells.Item(39, 2).Value = "Algorand"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(39, 4).Value = "0.2203"
$ws.Cells.Item(39, 5).Value = "  +3.10%  "
$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).Value = "9.067"
$ws.Cells.Item(40, 5).Value = "  +5.43%  "
$ws.Cells.Item(41, 2).Value = "TheSandbox"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(41, 4).Value = "0.6670"
$ws.Cells.Item(41, 5).Value = "  +3.55%  "
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "1.232"
$ws.Cells.Item(42, 5).Value = "  +2.36%  "
$ws.Cells.Item(43, 2).Value = "Frax"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(43, 4).Value = "1.001"
$ws.Cells.Item(43, 5).Value = "  +0.02%  "
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).Value = "13.75"
$ws.Cells.Item(44, 5).Value = "  +3.96%  "
$ws.Cells.Item(45, 2).Value = "Decentraland"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(45, 4).Value = "0.6186"
$ws.Cells.Item(45, 5).Value = "  +3.17%  "
$ws.Cells.Item(46, 2).Value = "NEARProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(46, 4).Value = "2.192"
$ws.Cells.Item(46, 5).Value = "  +0.12%  "
$ws.Cells.Item(47, 2).Value = "PancakeSwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(47, 4).Value = "3.659"
$ws.Cells.Item(47, 5).Value = "  +0.80%  "
$ws.Cells.Item(48, 2).Value = "EOS"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(48, 4).Value = "1.266"
$ws.Cells.Item(48, 5).Value = "  +5.33%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).Value = "124.68"
$ws.Cells.Item(49, 5).Value = "  +2.02%  "
$ws.Cells.Item(50, 2).Value = "Aave"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(50, 4).Value = "81.00"
$ws.Cells.Item(50, 5).Value = "  +3.51%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "0.06906"
$ws.Cells.Item(51, 5).Value = "  +1.80%  "
